# "Actualizo ICA y resultado fiscal marzo25"
# Update Expo-ICA, Impo-ICA and BC por zonas figures for March 2025.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("BC por zonas")
$values_BC_por_zonas = [ordered]@{
    "B2" = 3653.70565373
    "C2" = 5019.9640024500004
    "B3" = 2828.1663856
    "C3" = 4230.5198007099998
    "B4" = 2670.8467959899999
    "C4" = 532.69418169000005
    "B5" = 1509.1511670699999
    "C5" = 247.68493119999999
    "B6" = 614.63108608000005
    "C6" = 69.507858350000006
    "B7" = 1843.76028505
    "C7" = 2040.9377262800001
    "B8" = 1546.25451679
    "C8" = 1536.32525362
    "B9" = 1545.4638404499999
    "C9" = 2556.00661468
    "B10" = 940.04846015999999
    "C10" = 4354.9710388599997
    "B11" = 1053.4164926599999
    "C11" = 398.62831918000001
    "B12" = 1475.4662915700001
    "C12" = 1065.8297416
    "B13" = 738.60992347000001
    "C13" = 339.49022647999999
    "B14" = 304.34401658000002
    "C14" = 106.73525998
    "B15" = 1108.92552056
    "C15" = 99.697989109999995
    "B16" = 637.40449541999999
    "C16" = 107.71198874
    "B17" = 138.02896125999999
    "C17" = 39.985986320000002
    "B18" = 147.64052384999999
    "C18" = 105.53415789
    "B19" = 3168.58881703
    "C19" = 1300.8093220400001
}
foreach ($key in $values_BC_por_zonas.Keys) {
    $ws.Range($key).Value = $values_BC_por_zonas[$key]
}

$ws = $wb.Worksheets.Item("Expo-ICA")
$values_Expo_ICA = [ordered]@{
    "B2" = 18383.296137730002
    "B3" = 4603.2252014599999
    "B4" = 7.7104051499999997
    "B5" = 419.47604722
    "B6" = 37.858826399999998
    "B7" = 153.27298983
    "B8" = 133.74508589000001
    "B9" = 3322.9510893000001
    "B10" = 327.44592691000003
    "B11" = 12.53319778
    "B12" = 8.0573634900000002
    "B13" = 15.37204127
    "B14" = 109.50337096
    "B15" = 55.298857259999998
    "B16" = 6315.8819369000003
    "B17" = 819.80245958
    "B18" = 65.892660230000004
    "B19" = 314.29963931999998
    "B20" = 15.840670449999999
    "B21" = 17.40591512
    "B22" = 42.049211509999999
    "B23" = 169.68534578000001
    "B24" = 1912.2479316599999
    "B25" = 104.51563212000001
    "B26" = 217.09362447000001
    "B27" = 185.07986094
    "B28" = 2089.5512267200002
    "B29" = 71.859733210000002
    "B30" = 95.571712480000002
    "B31" = 23.0900985
    "B32" = 171.89621481
    "B33" = 4930.1587045400001
    "B34" = 1048.2152467599999
    "B35" = 231.95207525000001
    "B36" = 37.412688289999998
    "B37" = 3.6159849199999998
    "B38" = 107.58347292000001
    "B39" = 20.179592029999998
    "B40" = 1.9620223999999999
    "B41" = 27.16118264
    "B42" = 961.13869026999998
    "B43" = 401.49239516
    "B44" = 351.60597890999998
    "B45" = 1641.3083551499999
    "B46" = 6.0656587399999999
    "B47" = 90.465361099999996
    "B48" = 2534.03029483
    "B49" = 1461.49696713
    "B50" = 567.04802558999995
    "B51" = 24.110999629999998
    "B52" = 402.42192506000004
    "B53" = 78.952377420000005
}
foreach ($key in $values_Expo_ICA.Keys) {
    $ws.Range($key).Value = $values_Expo_ICA[$key]
}

$ws = $wb.Worksheets.Item("Impo-ICA")
$values_Impo_ICA = [ordered]@{
    "B2" = 17622.77106884
    "B3" = 3648.1198531199998
    "B4" = 2749.29462397
    "B5" = 313.04894249
    "B6" = 585.77628665999998
    "B7" = 5782.7834002400004
    "B8" = 651.50372373000005
    "B9" = 163.11292777
    "B10" = 261.46904121
    "B11" = 4250.73489393
    "B12" = 455.9628136
    "B13" = 662.01973711999995
    "B14" = 42.681152099999998
    "B15" = 619.33858501999998
    "B16" = 3867.1243676600002
    "B17" = 1469.7218238299999
    "B18" = 401.81034148999998
    "B19" = 1995.5922023400001
    "B20" = 2503.5203185300002
    "B21" = 218.74497688
    "B22" = 340.98948874000001
    "B23" = 228.08038034
    "B24" = 366.15075883999998
    "B25" = 695.39663073999998
    "B26" = 371.62422053
    "B27" = 282.53386246000002
    "B28" = 1018.39933536
    "B29" = 140.80405680999999
}
foreach ($key in $values_Impo_ICA.Keys) {
    $ws.Range($key).Value = $values_Impo_ICA[$key]
}


# Restore the selection on "BC por zonas" (B2:C19, activeCell B2) without
# leaving it as the active tab.
$wsBC = $wb.Worksheets.Item("BC por zonas")
$wsBC.Activate()
$wsBC.Range("B2:C19").Select()

# "Expo-ICA" ends up as the active/selected tab, cursor back at A1.
$wsExpo = $wb.Worksheets.Item("Expo-ICA")
$wsExpo.Activate()
$wsExpo.Range("A1").Select()
